# Updates the cryptocurrency price table (columns B-E, rows 2-51) to reflect
# the latest values scraped for this run, per the commit:
# "Updated cryptos list on Sun Feb 26 14:46:54 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several new Price (column D) values are plain decimal-looking strings
# (e.g. "304.36"). Left alone, Excel would auto-convert those into numeric
# values when assigned, but the sheet stores Price as literal text, so we
# mark those specific cells as Text ("@") before writing their value. Cells
# whose text already contains multiple "." separators (e.g. "23.263.57")
# are not auto-numeric and do not need this treatment.
$textForceAddrs = @("D4", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.263.57"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.606.50"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "304.36"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").Value = "0.3770"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "52.22"
$ws.Range("E8").Value = "  +5.26%  "

$ws.Range("D9").Value = "0.3632"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.08147"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "22.87"
$ws.Range("E13").Value = "  +1.64%  "

$ws.Range("D14").Value = "6.590"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "7.400"
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "1.606.33"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "94.08"
$ws.Range("E18").Value = "  +2.34%  "

$ws.Range("D19").Value = "0.06925"
$ws.Range("E19").Value = "  +1.47%  "

$ws.Range("D20").Value = "18.15"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "6.532"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "12.92"
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("D24").Value = "23.226.52"
$ws.Range("E24").Value = "  +0.88%  "

$ws.Range("D25").Value = "2.448"
$ws.Range("E25").Value = "  +3.89%  "

$ws.Range("D26").Value = "3.073"
$ws.Range("E26").Value = "  +9.51%  "

$ws.Range("D27").Value = "21.18"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").Value = "149.87"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").Value = "5.273"
$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("D30").Value = "135.59"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("D31").Value = "2.384"
$ws.Range("E31").Value = "  +2.26%  "

$ws.Range("D32").Value = "6.768"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").Value = "1.781.28"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").Value = "0.9654"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").Value = "0.07492"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("D36").Value = "0.02765"
$ws.Range("E36").Value = "  +2.40%  "

$ws.Range("D37").Value = "10.38"
$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").Value = "6.147"
$ws.Range("E39").Value = "  -1.61%  "

$ws.Range("D40").Value = "0.08798"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("D41").Value = "1.392"
$ws.Range("E41").Value = "  +2.19%  "

$ws.Range("D42").Value = "0.7092"
$ws.Range("E42").Value = "  +0.99%  "

$ws.Range("D43").Value = "12.51"
$ws.Range("E43").Value = "  +1.15%  "

$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("D45").Value = "0.6546"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").Value = "2.324"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "4.006"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").Value = "132.93"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "0.07948"
$ws.Range("E49").Value = "  +0.62%  "

$ws.Range("D50").Value = "1.209"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  -3.29%  "

